$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row stays the same (row 1). Rebuild data rows 2-10.

$ECs   = "ECs"
$FAPs  = "FAPs"
$sCs   = "sCs"
$Efnb2 = "Efnb2"
$Epha3 = "Epha3"

$rows = @(
  @{ A=$ECs;  D=$ECs;  E=3; F=1; G=20.94432133333333; H=62.832964;          I=0.7396577289668299; J=0.7396577289668298; K=1; L=0.3333333333333333; M=0.06698166666666668; N=0.200945;           O=0.003012576978541733; P=0.003012576978541732; Q=1.402885550108889;  R=12.62596995098;    S=0.002228275846285932; T=0.002228275846285932 },
  @{ A=$ECs;  D=$FAPs; E=3; F=1; G=20.94432133333333; H=62.832964;          I=0.7396577289668299; J=0.7396577289668298; K=3; L=1;                   M=22.03620333333333;  N=66.10861;           O=0.9911034191912899;   P=0.9911034191912899;   Q=461.5333235800045;  R=4153.79991222004;  S=0.7330773042102895;  T=0.7330773042102894 },
  @{ A=$ECs;  D=$sCs;  E=3; F=1; G=20.94432133333333; H=62.832964;          I=0.7396577289668299; J=0.7396577289668298; K=2; L=0.6666666666666666; M=0.130825;           N=0.392475;           O=0.005884003830168287; P=0.005884003830168287; Q=2.740040838433333;  R=24.6603675459;     S=0.004352148910254404; T=0.004352148910254403 },
  @{ A=$FAPs; D=$ECs;  E=3; F=1; G=2.327094666666667; H=6.981284;           I=0.08218235047311259;J=0.08218235047311258;K=1; L=0.3333333333333333; M=0.06698166666666668; N=0.200945;           O=0.003012576978541733; P=0.003012576978541732; Q=0.1558726792644445; R=1.40285411338;     S=0.0002475806570777473;T=0.0002475806570777472 },
  @{ A=$FAPs; D=$FAPs; E=3; F=1; G=2.327094666666667; H=6.981284;           I=0.08218235047311259;J=0.08218235047311258;K=3; L=1;                   M=22.03620333333333;  N=66.10861;           O=0.9911034191912899;   P=0.9911034191912899;   Q=51.28033125058223;  R=461.52298125524;   S=0.0814512085510788;   T=0.08145120855107879 },
  @{ A=$FAPs; D=$sCs;  E=3; F=1; G=2.327094666666667; H=6.981284;           I=0.08218235047311259;J=0.08218235047311258;K=2; L=0.6666666666666666; M=0.130825;           N=0.392475;           O=0.005884003830168287; P=0.005884003830168287; Q=0.3044421597666667; R=2.7399794379;      S=0.000483561264956027; T=0.0004835612649560269 },
  @{ A=$sCs;  D=$ECs;  E=3; F=1; G=5.044818;           H=15.134454;         I=0.1781599205600575; J=0.1781599205600575; K=1; L=0.3333333333333333; M=0.06698166666666668; N=0.200945;           O=0.003012576978541733; P=0.003012576978541732; Q=0.3379103176700001; R=3.041192859030001; S=0.0005367204751780532;T=0.0005367204751780531 },
  @{ A=$sCs;  D=$FAPs; E=3; F=1; G=5.044818;           H=15.134454;         I=0.1781599205600575; J=0.1781599205600575; K=3; L=1;                   M=22.03620333333333;  N=66.10861;           O=0.9911034191912899;   P=0.9911034191912899;   Q=111.16863522766;    R=1000.51771704894;  S=0.1765749064299216;   T=0.1765749064299216 },
  @{ A=$sCs;  D=$sCs;  E=3; F=1; G=5.044818;           H=15.134454;         I=0.1781599205600575; J=0.1781599205600575; K=2; L=0.6666666666666666; M=0.130825;           N=0.392475;           O=0.005884003830168287; P=0.005884003830168287; Q=0.65998831485;      R=5.939894833650001; S=0.001048293654957856; T=0.001048293654957856 }
)

$r = 2
foreach ($row in $rows) {
  $ws.Cells.Item($r, 1).Value  = $row.A
  $ws.Cells.Item($r, 2).Value  = $Efnb2
  $ws.Cells.Item($r, 3).Value  = $Epha3
  $ws.Cells.Item($r, 4).Value  = $row.D
  $ws.Cells.Item($r, 5).Value  = $row.E
  $ws.Cells.Item($r, 6).Value  = $row.F
  $ws.Cells.Item($r, 7).Value  = $row.G
  $ws.Cells.Item($r, 8).Value  = $row.H
  $ws.Cells.Item($r, 9).Value  = $row.I
  $ws.Cells.Item($r, 10).Value = $row.J
  $ws.Cells.Item($r, 11).Value = $row.K
  $ws.Cells.Item($r, 12).Value = $row.L
  $ws.Cells.Item($r, 13).Value = $row.M
  $ws.Cells.Item($r, 14).Value = $row.N
  $ws.Cells.Item($r, 15).Value = $row.O
  $ws.Cells.Item($r, 16).Value = $row.P
  $ws.Cells.Item($r, 17).Value = $row.Q
  $ws.Cells.Item($r, 18).Value = $row.R
  $ws.Cells.Item($r, 19).Value = $row.S
  $ws.Cells.Item($r, 20).Value = $row.T
  $r = $r + 1
}
